# Progress report updated by Suruchi:
# Her "animation and state management" update note was mis-entered in
# D10 (Sanskruti Nakhale's row) instead of D8 (Suruchi Shrey's own row).
# Move that note to the correct cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")   # same sheet as $wb.ActiveSheet here

# Value (getter) is unreliable in this host; use Value2 to read the text.
$note = $ws.Range("D10").Value2

$ws.Range("D8").Value = $note
$ws.Range("D10").ClearContents()

# Leave the cursor where Suruchi's edit ended up.
$ws.Range("D10").Select()
